$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep a Text format so numeric-looking strings
# (e.g. "227.18", "1.00", "0.0509") are stored as exact text, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value2 = "34.085.62"
$ws.Range("E2").Value2 = "  -0.01%  "
$ws.Range("D3").Value2 = "1.790.04"
$ws.Range("E3").Value2 = "  +0.32%  "
$ws.Range("D5").Value2 = "227.18"
$ws.Range("E5").Value2 = "  +0.65%  "
$ws.Range("D6").Value2 = "0.546"
$ws.Range("E6").Value2 = "  -0.63%  "
$ws.Range("E7").Value2 = "  +0.09%  "
$ws.Range("D8").Value2 = "32.16"
$ws.Range("E8").Value2 = "  -1.72%  "
$ws.Range("E9").Value2 = "  +2.91%  "
$ws.Range("E10").Value2 = "  -2.87%  "
$ws.Range("E11").Value2 = "  +0.32%  "
$ws.Range("D12").Value2 = "2.047.65"
$ws.Range("E12").Value2 = "  +0.29%  "
$ws.Range("D13").Value2 = "11.48"
$ws.Range("E13").Value2 = "  +4.46%  "
$ws.Range("D14").Value2 = "1.793.47"
$ws.Range("E14").Value2 = "  +0.42%  "
$ws.Range("D15").Value2 = "34.076.64"
$ws.Range("E15").Value2 = "  +0.15%  "
$ws.Range("D16").Value2 = "0.621"
$ws.Range("E16").Value2 = "  +0.14%  "
$ws.Range("E17").Value2 = "  +0.65%  "
$ws.Range("D18").Value2 = "67.82"
$ws.Range("E18").Value2 = "  +0.06%  "
$ws.Range("D19").Value2 = "245.41"
$ws.Range("E19").Value2 = "  +0.16%  "
$ws.Range("E20").Value2 = "  -0.96%  "
$ws.Range("B21").Value2 = "Avalanche"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value2 = "10.91"
$ws.Range("E21").Value2 = "  +1.25%  "
$ws.Range("B22").Value2 = "Dai"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value2 = "1.00"
$ws.Range("E22").Value2 = "  +0.05%  "
$ws.Range("E23").Value2 = "  +0.65%  "
$ws.Range("E24").Value2 = "  -2.38%  "
$ws.Range("D25").Value2 = "161.95"
$ws.Range("E25").Value2 = "  +1.11%  "
$ws.Range("D27").Value2 = "16.28"
$ws.Range("E27").Value2 = "  -0.28%  "
$ws.Range("D28").Value2 = "0.113"
$ws.Range("E28").Value2 = "  +0.93%  "
$ws.Range("E29").Value2 = "  +0.23%  "
$ws.Range("E30").Value2 = "  +1.14%  "
$ws.Range("E31").Value2 = "  +1.42%  "
$ws.Range("D32").Value2 = "3.66"
$ws.Range("E32").Value2 = "  +1.01%  "
$ws.Range("E33").Value2 = "  +3.01%  "
$ws.Range("E34").Value2 = "  +0.78%  "
$ws.Range("D35").Value2 = "1.437.17"
$ws.Range("E35").Value2 = "  +3.27%  "
$ws.Range("E36").Value2 = "  -0.54%  "
$ws.Range("E37").Value2 = "  +2.43%  "
$ws.Range("E38").Value2 = "  +6.79%  "
$ws.Range("D39").Value2 = "1.03"
$ws.Range("E39").Value2 = "  -1.13%  "
$ws.Range("D40").Value2 = "80.42"
$ws.Range("E40").Value2 = "  +2.94%  "
$ws.Range("E41").Value2 = "  +0.46%  "
$ws.Range("E42").Value2 = "  +0.69%  "
$ws.Range("E43").Value2 = "  +0.11%  "
$ws.Range("D44").Value2 = "13.34"
$ws.Range("E44").Value2 = "  +6.96%  "
$ws.Range("B45").Value2 = "BabyDogeCoin"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value2 = "0.0₆0139"
$ws.Range("E45").Value2 = "  -2.48%  "
$ws.Range("B46").Value2 = "Kaspa"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").Value2 = "0.0509"
$ws.Range("E46").Value2 = "  +2.43%  "
$ws.Range("B47").Value2 = "FraxShare"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value2 = "6.06"
$ws.Range("E47").Value2 = "  +3.90%  "
$ws.Range("E48").Value2 = "  -0.81%  "
$ws.Range("D49").Value2 = "107.49"
$ws.Range("E49").Value2 = "  -0.70%  "
$ws.Range("D50").Value2 = "1.949.29"
$ws.Range("E50").Value2 = "  +0.34%  "
$ws.Range("E51").Value2 = "  +0.08%  "
